$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "Reel N" -> "Image N" labels in column A (rows 6-8)
$ws.Range("A6").Value = "Image 1"
$ws.Range("A7").Value = "Image 2"
$ws.Range("A8").Value = "Image 3"

# Update D8 value 383 -> 384
$ws.Range("D8").Value = 384

# Update the active selection to I16
$ws.Range("I16").Select()
